$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 175, pushing existing rows 175-268 down to 177-270
$ws.Rows("175:176").Insert()

# Row 175: Primera, date 44452 (2021-09-13)
$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(175, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(175, 3).Value = "Coquimbo"
$ws.Cells.Item(175, 4).Value = 44452
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 100112037
$ws.Cells.Item(175, 7).Value = "Cebollín"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 320
$ws.Cells.Item(175, 11).Value = 3000
$ws.Cells.Item(175, 12).Value = 3500
$ws.Cells.Item(175, 13).Value = 3203
$ws.Cells.Item(175, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(175, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(175, 16).Value = 89
$ws.Cells.Item(175, 17).Value = 36
$ws.Cells.Item(175, 18).Value = "Hortaliza"

# Row 176: Segunda, date 44452
$ws.Cells.Item(176, 1).Value = 3
$ws.Cells.Item(176, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(176, 3).Value = "Coquimbo"
$ws.Cells.Item(176, 4).Value = 44452
$ws.Cells.Item(176, 5).Value = 5
$ws.Cells.Item(176, 6).Value = 100112037
$ws.Cells.Item(176, 7).Value = "Cebollín"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Segunda"
$ws.Cells.Item(176, 10).Value = 180
$ws.Cells.Item(176, 11).Value = 2500
$ws.Cells.Item(176, 12).Value = 2500
$ws.Cells.Item(176, 13).Value = 2500
$ws.Cells.Item(176, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(176, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(176, 16).Value = 69
$ws.Cells.Item(176, 17).Value = 36
$ws.Cells.Item(176, 18).Value = "Hortaliza"

# Apply date style (numFmtId 165 style) to new D cells by copying from D177 (previously D175)
$ws.Cells.Item(175, 4).NumberFormat = $ws.Cells.Item(177, 4).NumberFormat
$ws.Cells.Item(176, 4).NumberFormat = $ws.Cells.Item(177, 4).NumberFormat
